$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 69: date + commit description, continuing the log table.
$ws.Range("A69").Value = 41382
$ws.Range("B69").Value = "Wrote chapter platform model and first part of chapter execution model"

# Move the active selection to B70, matching where the next entry would be typed.
$ws.Range("B70").Select()
